# "Renouveau DonneeTest et Changement commentaire"
# Regenerate the test data on the "Resultat" sheet: refresh row 2 (summary
# values) and the column C weight list, extending it from 63 down to 75 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 : summary values (Nombres d'objets, Capacite du sac, etc.) ---
$ws.Range("A2").Value = 74
$ws.Range("B2").Value = 199
$ws.Range("C2").Value = 188
$ws.Range("E2").Value = 37
$ws.Range("F2").Value = 39
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 41
$ws.Range("I2").Value = 4

# --- Column C : the "Poids" data series, refreshed and extended to C75 ---
$weights = @(
    176, 85, 175, 31, 61, 192, 52, 122, 131, 73, 21, 173, 134, 51, 115,
    37, 125, 31, 143, 98, 52, 196, 69, 60, 162, 152, 89, 82, 129, 39,
    190, 124, 85, 164, 152, 121, 150, 76, 1, 110, 105, 65, 74, 11, 105,
    119, 38, 69, 164, 143, 162, 102, 14, 65, 54, 117, 60, 41, 115, 143,
    71, 39, 197, 53, 0, 162, 74, 148, 11, 19, 90, 158, 144
)

$startRow = 3
for ($i = 0; $i -lt $weights.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 3).Value = $weights[$i]
}
